# Trade #109 closed at 2026-02-17 09:19:12 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.12   # Current Capital
$summary.Range("B4").Value = 0.13      # Total P&L $
$summary.Range("B6").Value = 109       # Total Trades
$summary.Range("B7").Value = 47        # Winning Trades
$summary.Range("B9").Value = 43.12     # Win Rate %

# --- Strategy Status sheet (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.12     # Capital
$status.Range("D4").Value = 109        # Trades
$status.Range("E4").Value = 0.13       # P&L $
$status.Range("F4").Value = 0.12       # P&L %
$status.Range("G4").Value = 43.12      # Win Rate %

# --- New trade row (#109) appended to "All Trades" and "MarketMaking" sheets ---
$newRow = 110

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item($newRow, 1).Value = 109                 # A Trade #

    # Date column: force text so Excel doesn't auto-convert it to a date serial
    $ws.Cells.Item($newRow, 2).NumberFormat = "@"
    $ws.Cells.Item($newRow, 2).Value = "2026-02-17"        # B Date

    $ws.Cells.Item($newRow, 3).Value = "09:19:05"          # C Time
    $ws.Cells.Item($newRow, 4).Value = "MarketMaking"      # D Strategy
    $ws.Cells.Item($newRow, 5).Value = "UP"                # E Side
    $ws.Cells.Item($newRow, 6).Value = 0.19                # F Entry Price
    $ws.Cells.Item($newRow, 7).Value = 0.21                # G Exit Price
    $ws.Cells.Item($newRow, 8).Value = "CLOSED"            # H Status
    $ws.Cells.Item($newRow, 9).Value = 10.5263             # I P&L %
    $ws.Cells.Item($newRow, 10).Value = 0.02               # J P&L $
    $ws.Cells.Item($newRow, 11).Value = 100.12             # K Capital After
    $ws.Cells.Item($newRow, 12).Value = 0                  # L Entry Slippage (bps)
    $ws.Cells.Item($newRow, 13).Value = 0                  # M Exit Slippage (bps)
    $ws.Cells.Item($newRow, 14).Value = 0.6                # N Confidence
    $ws.Cells.Item($newRow, 15).Value = "Normal spread capture: 19600 bps"  # O Entry Reason
    $ws.Cells.Item($newRow, 16).Value = "early_exit"       # P Exit Reason
    $ws.Cells.Item($newRow, 17).Value = 0.14               # Q Duration (min)
}
